# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp banner
$ws.Range("A1").Value = "Datos actualizados a 13 de Septiembre de 2020 a las 09:35"

# Row 55: Singapur (grew in cases) - rank unchanged
$ws.Range("B55").Value = 57406
$ws.Range("C55").Value = 49
$ws.Range("E55").Value = 680

# Row 62: Armenia - rank unchanged
$ws.Range("B62").Value = 45862
$ws.Range("C62").Value = 187
$ws.Range("D62").Value = 41659
$ws.Range("E62").Value = 3287
$ws.Range("G62").Value = 5
$ws.Range("H62").Value = 916

# Rows 92/93: Hungria overtakes Noruega in ranking, so they swap places.
# Row 92 becomes Hungria with its updated (grown) figures.
$ws.Range("A92").Value = "Hungria"
$ws.Range("B92").Value = 12309
$ws.Range("C92").Value = 484
$ws.Range("D92").Value = 4069
$ws.Range("E92").Value = 7603
$ws.Range("G92").Value = 4
$ws.Range("H92").Value = 637

# Row 93 becomes Noruega, keeping its previous (now older) figures.
$ws.Range("A93").Value = "Noruega"
$ws.Range("B93").Value = 12079
$ws.Range("C93").Value = 0
$ws.Range("D93").Value = 10371
$ws.Range("E93").Value = 1443
$ws.Range("G93").Value = 0
$ws.Range("H93").Value = 265

# Rows 149/150/151: Georgia overtakes Islandia and Sierra Leona, so the
# three rows rotate: Georgia moves up to row 149 with updated figures,
# Islandia drops to row 150, Sierra Leona drops to row 151 (each keeping
# its previous figures).
$ws.Range("A149").Value = "Georgia"
$ws.Range("B149").Value = 2227
$ws.Range("C149").Value = 152
$ws.Range("D149").Value = 1369
$ws.Range("E149").Value = 839
$ws.Range("H149").Value = 19

$ws.Range("A150").Value = "Islandia"
$ws.Range("B150").Value = 2162
$ws.Range("C150").Value = 0
$ws.Range("D150").Value = 2085
$ws.Range("E150").Value = 67
$ws.Range("H150").Value = 10

$ws.Range("A151").Value = "Sierra Leona"
$ws.Range("B151").Value = 2096
$ws.Range("C151").Value = 0
$ws.Range("D151").Value = 1634
$ws.Range("E151").Value = 390
$ws.Range("H151").Value = 72

# Row 159: Letonia - rank unchanged
$ws.Range("B159").Value = 1474
$ws.Range("C159").Value = 10
$ws.Range("E159").Value = 191
